$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.028.40'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').Value = '2.310.21'
$ws.Range('E3').Value = '  -3.04%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.79%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.510'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.509'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '51.31'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0797'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.03%  '
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.80'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.96%  '
$ws.Range('D15').Value = '2.660.56'
$ws.Range('E15').Value = '  -3.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').Value = '2.293.15'
$ws.Range('E17').Value = '  -3.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.801'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.46%  '
$ws.Range('D19').Value = '42.885.27'
$ws.Range('E19').Value = '  -1.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.79'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('D21').Value = '0.0₃0903'
$ws.Range('E21').Value = '  -1.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.52'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '236.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.98'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.61%  '
$ws.Range('E26').Value = '  -3.71%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.95'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.78'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '164.94'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.05'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.16%  '
$ws.Range('E35').Value = '  -4.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0708'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.51'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -10.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.87'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.54%  '
$ws.Range('E40').Value = '  -7.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.102'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.09%  '
$ws.Range('E42').Value = '  -2.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.41'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.56%  '
$ws.Range('D44').Value = '1.974.60'
$ws.Range('E44').Value = '  -3.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0286'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '18.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.87'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.91'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.84%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.81'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '54.26'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.50%  '
$ws.Range('D51').Value = '2.529.87'
$ws.Range('E51').Value = '  -2.76%  '
